$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Force text format on the cells we will touch to avoid Excel auto-converting
# numeric-looking strings (e.g. "235.32", "1.001") into actual numbers.
$cells = @("D2","E2","D3","E3","E4","D5","E5","D6","E6","D7","E7","D8","E8","B9","C9","D9","E9","B10","C10","D10","E10","B11","C11","D11","E11","B12","C12","D12","E12","B13","C13","D13","E13","B14","C14","D14","E14","B15","C15","D15","E15","B16","C16","D16","E16","B17","C17","B18","C18","D18","E18","B19","C19","D19","E19","B20","C20","D20","E20","B21","C21","D21","E21","B22","C22","D22","E22","B23","C23","D23","E23","B24","C24","D24","E24","B25","C25","D25","E25","B26","C26","D26","E26","B27","C27","D27","E27","D28","E28","B29","C29","D29","E29","B30","C30","D30","E30","B31","C31","D31","E31","B32","C32","D32","E32","B33","C33","D33","E33","B34","C34","D34","E34","B35","C35","D35","E35","B36","C36","D36","E36","B37","C37","D37","E37","B38","C38","D38","E38","B39","C39","D39","E39","B40","C40","D40","E40","B41","C41","D41","E41","B42","C42","D42","E42","B43","C43","D43","E43","B44","C44","D44","E44","B45","C45","D45","E45","B46","C46","D46","E46","B47","C47","D47","E47","B48","C48","D48","E48","B49","C49","D49","E49","B50","C50","D50","E50","B51","C51","D51","E51")
foreach ($c in $cells) { $ws.Range($c).NumberFormat = "@" }

$ws.Range("D2").Value = '25.926.10'
$ws.Range("E2").Value = '  +0.17%  '
$ws.Range("D3").Value = '1.753.46'
$ws.Range("E3").Value = '  -0.07%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '235.32'
$ws.Range("E5").Value = '  -1.67%  '
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  +0.11%  '
$ws.Range("D7").Value = '0.5205'
$ws.Range("E7").Value = '  +2.17%  '
$ws.Range("D8").Value = '0.2724'
$ws.Range("E8").Value = '  -1.51%  '
$ws.Range("B9").Value = 'Dogecoin'
$ws.Range("C9").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D9").Value = '0.06154'
$ws.Range("E9").Value = '  -0.68%  '
$ws.Range("B10").Value = 'WrappedEther'
$ws.Range("C10").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D10").Value = '1.764.05'
$ws.Range("E10").Value = '  +0.58%  '
$ws.Range("B11").Value = 'TRON'
$ws.Range("C11").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D11").Value = '0.07030'
$ws.Range("E11").Value = '  +0.75%  '
$ws.Range("B12").Value = 'Solana'
$ws.Range("C12").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D12").Value = '15.53'
$ws.Range("E12").Value = '  -1.52%  '
$ws.Range("B13").Value = 'Polygon'
$ws.Range("C13").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D13").Value = '0.6348'
$ws.Range("E13").Value = '  +3.57%  '
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").Value = '4.516'
$ws.Range("E14").Value = '  -0.44%  '
$ws.Range("B15").Value = 'Litecoin'
$ws.Range("C15").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D15").Value = '77.57'
$ws.Range("E15").Value = '  +0.22%  '
$ws.Range("B16").Value = 'BinanceUSD'
$ws.Range("C16").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D16").Value = '1.001'
$ws.Range("E16").Value = '  +0.08%  '
$ws.Range("B17").Value = 'Dai'
$ws.Range("C17").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").Value = '25.937.02'
$ws.Range("E18").Value = '  +0.16%  '
$ws.Range("B19").Value = 'Avalanche'
$ws.Range("C19").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D19").Value = '11.56'
$ws.Range("E19").Value = '  -1.19%  '
$ws.Range("B20").Value = 'ShibaInu'
$ws.Range("C20").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D20").Value = '0.000006644'
$ws.Range("E20").Value = '  -3.81%  '
$ws.Range("B21").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C21").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D21").Value = '1.989.51'
$ws.Range("E21").Value = '  +0.83%  '
$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").Value = '4.053'
$ws.Range("E22").Value = '  -0.73%  '
$ws.Range("B23").Value = 'Cosmos'
$ws.Range("C23").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D23").Value = '8.429'
$ws.Range("E23").Value = '  +2.12%  '
$ws.Range("B24").Value = 'Chainlink'
$ws.Range("C24").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D24").Value = '5.161'
$ws.Range("E24").Value = '  -1.94%  '
$ws.Range("B25").Value = 'Monero'
$ws.Range("C25").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D25").Value = '138.71'
$ws.Range("E25").Value = '  +0.41%  '
$ws.Range("B26").Value = 'Toncoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D26").Value = '1.502'
$ws.Range("E26").Value = '  +0.60%  '
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").Value = '15.07'
$ws.Range("E27").Value = '  -0.13%  '
$ws.Range("D28").Value = '1.819'
$ws.Range("E28").Value = '  -0.06%  '
$ws.Range("B29").Value = 'BitcoinCash'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D29").Value = '103.01'
$ws.Range("E29").Value = '  -0.92%  '
$ws.Range("B30").Value = 'Stellar'
$ws.Range("C30").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D30").Value = '0.08335'
$ws.Range("E30").Value = '  +1.02%  '
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").Value = '3.655'
$ws.Range("E31").Value = '  -1.26%  '
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").Value = '3.419'
$ws.Range("E32").Value = '  -2.24%  '
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").Value = '0.04421'
$ws.Range("E33").Value = '  -2.97%  '
$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D34").Value = '2.630'
$ws.Range("E34").Value = '  -0.58%  '
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").Value = '0.9918'
$ws.Range("E35").Value = '  -0.20%  '
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").Value = '0.6023'
$ws.Range("E36").Value = '  -1.50%  '
$ws.Range("B37").Value = 'MXToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D37").Value = '2.710'
$ws.Range("E37").Value = '  +0.39%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = '0.01579'
$ws.Range("E38").Value = '  +1.11%  '
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").Value = '1.947'
$ws.Range("E39").Value = '  +2.58%  '
$ws.Range("B40").Value = 'PaxDollar'
$ws.Range("C40").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D40").Value = '1.001'
$ws.Range("E40").Value = '  +0.05%  '
$ws.Range("B41").Value = 'Quant'
$ws.Range("C41").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D41").Value = '102.42'
$ws.Range("E41").Value = '  -1.31%  '
$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D42").Value = '0.3843'
$ws.Range("E42").Value = '  -0.51%  '
$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").Value = '0.7372'
$ws.Range("E43").Value = '  -0.60%  '
$ws.Range("B44").Value = 'FraxShare'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D44").Value = '4.889'
$ws.Range("E44").Value = '  -1.81%  '
$ws.Range("B45").Value = 'Cronos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D45").Value = '0.05509'
$ws.Range("E45").Value = '  +1.41%  '
$ws.Range("B46").Value = 'Aptos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D46").Value = '6.263'
$ws.Range("E46").Value = '  +3.93%  '
$ws.Range("B47").Value = 'Algorand'
$ws.Range("C47").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D47").Value = '0.1106'
$ws.Range("E47").Value = '  -0.85%  '
$ws.Range("B48").Value = 'Elrond'
$ws.Range("C48").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D48").Value = '30.02'
$ws.Range("E48").Value = '  -0.15%  '
$ws.Range("B49").Value = 'USDD'
$ws.Range("C49").Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range("D49").Value = '1.003'
$ws.Range("E49").Value = '  +0.51%  '
$ws.Range("B50").Value = 'Aave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D50").Value = '52.17'
$ws.Range("E50").Value = '  -1.12%  '
$ws.Range("B51").Value = 'Decentraland'
$ws.Range("C51").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D51").Value = '0.3408'
$ws.Range("E51").Value = '  -1.66%  '
